$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "KeyValuePairs" table (Tabelle2) currently covers A1:C47 - two new
# localization rows are being appended for the service-worker update modal,
# so grow the table (and with it the AutoFilter range) down to row 49.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C49"))

# Copy the formatting (alternating-row fill) of the last existing data row
# onto the two new rows before filling in their values.
$ws.Range("A46:C46").Copy()
$ws.Range("A48:C49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new localization keys / strings, column by column (Key first,
# then the English strings, then the German strings) for the two new rows.
$ws.Range("A48").Value = "service_worker-update_headline"
$ws.Range("A49").Value = "service_worker-update_confirm_btn_txt"

$ws.Range("B48").Value = "Update Available"
$ws.Range("C48").Value = "Update verfügbar"

$ws.Range("B49").Value = "Update Now & Refresh"
$ws.Range("C49").Value = "Update installieren"

# The longer key in column A pushes the "best fit" column width out a bit.
$ws.Columns.Item(1).ColumnWidth = 31.84

# Leave the view scrolled down to / focused on the newly added rows.
$win = $excel.Windows.Item(1)
$win.ScrollRow = 31
$ws.Range("D44").Select()
